# Insert a new weekly price record as row 385, shifting all subsequent
# records (old rows 385-441) down by one row (new rows 386-442).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(385).Insert()

$ws.Range("A385").Value = 4
$ws.Range("B385").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C385").Value = "Los Lagos"
$ws.Range("D385").Value = 45077
$ws.Range("E385").Value = 10
$ws.Range("F385").Value = 100112043
$ws.Range("G385").Value = "Pepino ensalada"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 70
$ws.Range("K385").Value = 16000
$ws.Range("L385").Value = 16000
$ws.Range("M385").Value = 16000
$ws.Range("N385").Value = "`$/caja 60 unidades"
$ws.Range("O385").Value = "Región de Arica y Parinacota"
$ws.Range("P385").Value = 267
$ws.Range("Q385").Value = 60
$ws.Range("R385").Value = "Hortaliza"
